$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the auto "datetimeFigureOut" date placeholder that is cached on
#    the slide master and on every slide layout (4/18/2020 -> 4/21/2020).
# ---------------------------------------------------------------------------
function Update-DatePlaceholder {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)

        $isPlaceholder = $false
        try { $isPlaceholder = ($sh.Type -eq 14) } catch { $isPlaceholder = $false }
        if (-not $isPlaceholder) { continue }

        $phType = -1
        try { $phType = $sh.PlaceholderFormat.Type } catch { $phType = -1 }

        if ($phType -eq 16) {
            # ppPlaceholderDate
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "4/21/2020"

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes $newDate

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholder $layout.Shapes $newDate
}

# ---------------------------------------------------------------------------
# 2) Reposition / resize the two result pictures on slide 1.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)

$picGrid = $slide.Shapes.Item(20)     # "Picture 2" - grid based map picture
$picGrid.Left   = 17.79984251968504
$picGrid.Top    = 1726.4837007874016
$picGrid.Width  = 1091.9153543307086
$picGrid.Height = 683.9993700787402

$picDimRed = $slide.Shapes.Item(27)   # "Picture 8" - dimensionality reduction picture
$picDimRed.Left   = 2313.605118110236
$picDimRed.Top    = 826.3567716535433
$picDimRed.Width  = 1102.1433070866142
$picDimRed.Height = 629.5152755905511
